$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was updated from 45192
# (2023-09-23) to 45202 (2023-10-03) for every data row (rows 2 through 110).
for ($r = 2; $r -le 110; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value2 = 45202
    }
}
